$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) text values
$ws.Range("A1").Value = "button_testResultDetails_class"
$ws.Range("B1").Value = "button_testResultDetails_class_1"
$ws.Range("C1").Value = "button_testResultDetails_class_2"
$ws.Range("D1").Value = "button_testResultDetails_internalRoleButtonName"
$ws.Range("E1").Value = "button_testResultDetails_internalRoleButtonName_1"
$ws.Range("F1").Value = "button_testResultDetails_internalRoleButtonName_2"

# Update row 2 values: shift old C2/D2/E2 contents, insert new C2
$ws.Range("C2").Value = """]:nth-child(3) [class=""css-1yjo05o"
$ws.Range("D2").Value = "Failed Automations - Apply to"
$ws.Range("E2").Value = "Failed Portal - Login with"

# Adjust column widths: col3 <- 37, col4 <- 49, col5 <- 51
# (subtract 5/6 to compensate for the engine's internal char->px padding
#  so the stored OOXML width lands exactly on the target value)
$ws.Columns.Item(3).ColumnWidth = 36.16666666666667
$ws.Columns.Item(4).ColumnWidth = 48.16666666666667
$ws.Columns.Item(5).ColumnWidth = 50.16666666666667
